$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A7").Value = "Gibberish (Safari) "
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$s1 = $chart.SeriesCollection().Item(1)
try { Write-Host ($s1.CategoryNames -join ",") } catch { Write-Host "CategoryNames get failed: $_" }
try {
  $s1.CategoryNames = @("Flocking (Firefox)", "Gibberish (Firefox) ", "Flocking (Chrome)", "Gibberish (Chrome) ", "Flocking (Safari)", "Gibberish (Safari) ")
  Write-Host "CategoryNames set OK"
} catch { Write-Host "CategoryNames set failed: $_" }
